$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original text formatting for Price (D) and Hora (G) columns
# so updated values remain stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Update Hora (G) column: all rows 2-51 change from 9 to 10
$ws.Range("G2:G51").Value = "10"

# Update Price (D) column for rows with refreshed price data
$ws.Range("D2").Value = "243.58"
$ws.Range("D3").Value = "23.06"
$ws.Range("D4").Value = "5.405"
$ws.Range("D5").Value = "0.05961"
$ws.Range("D6").Value = "3.428"
$ws.Range("D7").Value = "6.502"
$ws.Range("D8").Value = "0.8111"
$ws.Range("D9").Value = "0.9214"
$ws.Range("D10").Value = "0.1433"
$ws.Range("D11").Value = "0.07430"
$ws.Range("D12").Value = "0.03261"
$ws.Range("D13").Value = "0.03074"
$ws.Range("D14").Value = "0.09356"
$ws.Range("D15").Value = "3.855"
$ws.Range("D16").Value = "0.001563"
$ws.Range("D17").Value = "0.04716"
$ws.Range("D19").Value = "0.005983"
$ws.Range("D20").Value = "0.001260"
$ws.Range("D21").Value = "0.004797"
$ws.Range("D22").Value = "0.00007995"
$ws.Range("D23").Value = "3.578"
$ws.Range("D25").Value = "0.3241"
$ws.Range("D27").Value = "0.0002338"
$ws.Range("D41").Value = "0.006366"
$ws.Range("D42").Value = "0.1076"
$ws.Range("D43").Value = "0.002538"
$ws.Range("D44").Value = "0.008904"
$ws.Range("D45").Value = "0.00005156"
$ws.Range("D47").Value = "0.6847"
$ws.Range("D48").Value = "0.002144"
